# Weekly fruit/vegetable price refresh: the daily-logic subset sheet gets
# re-shuffled into a new row order (each row keeps its original field values,
# but rows are re-sequenced). Re-derive this by snapshotting every data row
# (rows 2-29, columns A-R) then writing the rows back out in the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 29
$lastCol = 18   # column R

# Snapshot every existing data row before we overwrite anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $row = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $row += ,($ws.Cells.Item($r, $c).Value2)
    }
    $snapshot[$r] = $row
}

# new row number -> old row number that now supplies its data.
$mapping = @{
    2  = 27
    3  = 21
    4  = 20
    5  = 19
    6  = 15
    7  = 4
    8  = 23
    9  = 24
    10 = 17
    11 = 18
    12 = 28
    13 = 22
    14 = 12
    15 = 16
    16 = 3
    17 = 26
    18 = 11
    19 = 25
    20 = 14
    21 = 8
    22 = 9
    23 = 5
    24 = 7
    25 = 10
    26 = 2
    27 = 6
    28 = 29
    29 = 13
}

foreach ($newRow in ($mapping.Keys | Sort-Object)) {
    $oldRow = $mapping[$newRow]
    $data = $snapshot[$oldRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value2 = $data[$c - 1]
    }
}
